# The "Specification" section body paragraph is currently empty (the
# empty paragraph right after the "Specification" heading, before the
# "Programming" heading). Fill it in with the model-specification text,
# give the run 12pt (sz/szCs 24 half-points) sizing to match the rest of
# the document's body text, and add 12pt spacing before/after the
# paragraph.

$d = $word.ActiveDocument
$p = $d.Paragraphs(6)
$r = $p.Range

$r.Font.Size = 12
$r.Font.SizeBi = 12
$r.Text = "To examine the fairness of ""one person, one vote,"" we need to analyze population distribution and districting. In our model, we set the desired district population at 750,000 with six districts. We then calculate the ideal population by dividing the state's total population by the number of districts. The objective function calculates the absolute difference between the projected population of each district and this ideal population. Since we aim to minimize deviation from the ideal population regardless of whether it's an over- or underrepresentation, the objective function seeks to minimize the total absolute population deviation across all districts. The optimization problem then seeks to find the optimal assignment of counties to districts (represented by variables) that minimizes this deviation, potentially incorporating other relevant objectives."

$p.Format.SpaceBefore = 12
$p.Format.SpaceAfter = 12
